$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "2025-11-03T00:10:22.618693"
    3  = "2025-11-03T00:10:22.618693"
    4  = "2025-11-03T00:10:22.618693"
    5  = "2025-11-03T00:10:22.618693"
    6  = "2025-11-03T00:10:22.618693"
    7  = "2025-11-03T00:10:22.618693"
    8  = "2025-11-03T00:10:22.618693"
    9  = "2025-11-03T00:10:22.618693"
    10 = "2025-11-03T00:10:22.618693"
    11 = "2025-11-03T00:10:22.619693"
    12 = "2025-11-03T00:10:22.619693"
    13 = "2025-11-03T00:10:22.619693"
    14 = "2025-11-03T00:10:22.619693"
    15 = "2025-11-03T00:10:22.619693"
    16 = "2025-11-03T00:10:22.619693"
    17 = "2025-11-03T00:10:22.619693"
    18 = "2025-11-03T00:10:22.619693"
    19 = "2025-11-03T00:10:22.619693"
    20 = "2025-11-03T00:10:22.619693"
    21 = "2025-11-03T00:10:22.619693"
    22 = "2025-11-03T00:10:22.620692"
    23 = "2025-11-03T00:10:22.620692"
    24 = "2025-11-03T00:10:22.620692"
    25 = "2025-11-03T00:10:22.620692"
    26 = "2025-11-03T00:10:22.620692"
    27 = "2025-11-03T00:10:22.620692"
    28 = "2025-11-03T00:10:22.620692"
    29 = "2025-11-03T00:10:22.620692"
}

foreach ($row in $values.Keys) {
    $ws.Range("Z$row").Value = $values[$row]
}
